$d = $word.ActiveDocument

$pairs = @(
    @("2024-06-04 Tuesday", "2024-06-05 Wednesday"),
    @("300×3=900", "287×6=1722"),
    @("710×5=3550", "523×8=4184"),
    @("423×2=846", "318×8=2544"),
    @("108×9=972", "947×9=8523"),
    @("971×5=4855", "200×3=600"),
    @("238×6=1428", "693×6=4158"),
    @("114×6=684", "113×9=1017"),
    @("973×3=2919", "670×8=5360"),
    @("192×4=768", "607×6=3642"),
    @("241×3=723", "686×5=3430"),
    @("769×7=5383", "852×9=7668"),
    @("348×8=2784", "559×4=2236"),
    @("415×7=2905", "730×9=6570"),
    @("688×7=4816", "554×5=2770"),
    @("793×5=3965", "310×9=2790"),
    @("978×9=8802", "133×8=1064"),
    @("458×7=3206", "450×7=3150"),
    @("456×8=3648", "129×5=645"),
    @("585×6=3510", "509×4=2036"),
    @("582×4=2328", "736×4=2944"),
    @("766×4=3064", "412×3=1236"),
    @("908×3=2724", "425×6=2550"),
    @("147×4=588", "916×7=6412"),
    @("552×6=3312", "638×5=3190"),
    @("654×3=1962", "102×3=306")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
